# The workbook used to be produced like a raw CSV dump: the scraped
# laptop title/link pairs were written straight into columns C:D with no
# header row. This turns that into a proper Excel table living in A:B,
# with a "Titlu laptop" / "Link laptop" header on row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the block of data as it currently exists (columns C:D).
$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$lastRow = $firstRow + $rowCount - 1
$lastCol = $firstCol + $colCount - 1

$dataRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$data = $dataRange.Value()

# Remove the old data block entirely.
$dataRange.ClearContents()

# Write the new header row in A1:B1.
$ws.Cells.Item(1, 1).Value = "Titlu laptop"
$ws.Cells.Item(1, 2).Value = "Link laptop"

# Re-write the captured data into A:B, shifted down one row to leave room
# for the header.
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $ws.Cells.Item($r + 1, $c).Value = $data[$r, $c]
    }
}
